# Removed Test Case Inter-Dependency
#
# The loan product name previously embedded a hard-coded scenario suffix
# ("...-ONTIME-PER") and a numeric short name (4290) that collided with
# other test cases. Update both sheets to use an independent product
# name/short name, then leave the "ProductLoanOutput" sheet as the active
# tab with the input sheet's selection reset to the top of the form.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("ProductLoanInput")
$ws2 = $wb.Worksheets.Item("ProductLoanOutput")

$newProductName = "4290-MS-EI-DB-SAR-REC-RNI-FEE-FFConMONTHLYonDAY25-FIFC-1-FFROP-DAILY-FIFR-1-MD-TR-1-1st"
$newShortName   = "429t"

# productname (row 1) on both the input and output sheets
$ws1.Range("B1").Value = $newProductName
$ws2.Range("B1").Value = $newProductName

# shortname (row 2) on the input sheet -- now a free-form code instead of
# the raw numeric id that other test cases also happened to use
$ws1.Range("B2").Value = $newShortName

# Reset the input sheet's selection (it had been left on B17) and make
# the output sheet the active tab, matching a natural "review the output"
# flow after filling in the input sheet.
$ws1.Range("B3").Select() | Out-Null
$ws2.Activate() | Out-Null
